$d = $word.ActiveDocument

$replacements = @(
    @{old="2025-08-23 Saturday"; new="2025-08-24 Sunday"},
    @{old="36×91=3276"; new="96×52=4992"},
    @{old="94×14=1316"; new="93×22=2046"},
    @{old="84×83=6972"; new="46×82=3772"},
    @{old="19×53=1007"; new="21×37=777"},
    @{old="55×92=5060"; new="65×14=910"},
    @{old="11×90=990"; new="65×60=3900"},
    @{old="36×88=3168"; new="55×47=2585"},
    @{old="34×16=544"; new="91×24=2184"},
    @{old="93×61=5673"; new="19×18=342"},
    @{old="53×61=3233"; new="51×73=3723"},
    @{old="74×86=6364"; new="39×70=2730"},
    @{old="85×96=8160"; new="32×36=1152"},
    @{old="29×54=1566"; new="60×52=3120"},
    @{old="18×60=1080"; new="68×35=2380"},
    @{old="95×76=7220"; new="23×65=1495"},
    @{old="74×39=2886"; new="37×65=2405"},
    @{old="18×76=1368"; new="44×58=2552"},
    @{old="61×34=2074"; new="87×39=3393"},
    @{old="58×27=1566"; new="87×90=7830"},
    @{old="27×84=2268"; new="84×73=6132"},
    @{old="91×49=4459"; new="42×16=672"},
    @{old="96×32=3072"; new="86×36=3096"},
    @{old="56×98=5488"; new="19×66=1254"},
    @{old="58×91=5278"; new="81×28=2268"},
    @{old="21×24=504"; new="62×61=3782"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $r.new, 2)
}
